# "added knesset speakers function"
# The old "AutoAgg" placeholder rectangle (shape id 4, name "מלבן 3") on the
# title slide is no longer needed and is removed.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 4) {
        $target = $sh
        break
    }
}

if ($target -eq $null) {
    # Fall back to locating it by its leftover "AutoAgg" run text, in case
    # shape ids ever shift.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "AutoAgg") {
                $target = $sh
                break
            }
        }
    }
}

if ($target -ne $null) {
    $target.Delete()
}
